$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.263.14"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "1.824.79"
$ws.Range("E3").Value = "  -0.85%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("D5").Value = "'236.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.44%  "
$ws.Range("D6").Value = "'0.5994"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.32%  "
$ws.Range("D7").Value = "'1.005"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.47%  "
$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").Value = "'0.07117"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.96%  "
$ws.Range("B9").Value = "Solana"
$ws.Range("C9").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D9").Value = "'24.17"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.17%  "
$ws.Range("B10").Value = "Cardano"
$ws.Range("C10").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D10").Value = "'0.2801"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.11%  "
$ws.Range("D11").Value = "'0.07687"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.59%  "
$ws.Range("D12").Value = "1.844.68"
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").Value = "'4.777"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.95%  "
$ws.Range("D14").Value = "'0.6421"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.01%  "
$ws.Range("D15").Value = "'0.000009767"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.37%  "
$ws.Range("D16").Value = "'79.50"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.99%  "
$ws.Range("B17").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C17").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D17").Value = "2.038.75"
$ws.Range("E17").Value = "  -2.30%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Value = "'6.064"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.74%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "29.260.51"
$ws.Range("E19").Value = "  -0.19%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'229.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.53%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "'1.003"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.36%  "
$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D22").Value = "'11.78"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.02%  "
$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").Value = "'7.013"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.28%  "
$ws.Range("B24").Value = "BinanceUSD"
$ws.Range("C24").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D24").Value = "'1.002"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "'156.60"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.40%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "'8.110"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.08%  "
$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").Value = "'0.1274"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.50%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'16.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.33%  "
$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").Value = "'0.06802"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.63%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'1.472"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.58%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.462"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.55%  "
$ws.Range("D32").Value = "'3.809"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.20%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'3.758"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.35%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "'1.135"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Value = "'1.713"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.44%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.6621"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.16%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "'2.573"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.22%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").Value = "'2.766"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.99%  "
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "1.221.65"
$ws.Range("E39").Value = "  -1.68%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.01769"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.46%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'6.579"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.52%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'0.9278"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.59%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "'1.005"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.59%  "
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "1.961.63"
$ws.Range("E44").Value = "  -1.07%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "'99.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.41%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'63.44"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.05%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "'0.00000000118"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.03%  "
$ws.Range("D48").Value = "'1.633"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.35%  "
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").Value = "'6.595"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.16%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.05601"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.28%  "
$ws.Range("D51").Value = "'0.1087"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.37%  "
